$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G holds numeric-looking text ("0.0") for detail rows; force text
# formatting on those cells so Excel does not coerce them to the number 0.
$ws.Range("G2:G9,G11:G18,G21:G31,G33:G43").NumberFormat = "@"

# Full target table for the "manuallyCalcDif" sheet (A1:K44).
$table = @(
    @("Subscriber", "Application", "API", "Version", "Operation", "Charging Plan", "Count", "Usage Charge", "Tax", "Credit", "Grand Total"),
    @("AutInvSP3", "AutInvSP3_AutInvSP4_app4", "smsmessaging", "v1", "Send SMS", "SM1", "0.0", "Rs 0.0", "Rs 0.0", "Rs 0.0", ""),
    @("", "", "", "", "Retrive SMS", "SM2", "0.0", "Rs 0.0", "Rs 0.0", "Rs 0.0", ""),
    @("", "", "", "", "Query SMS Delivery", "SM2", "0.0", "Rs 0.0", "Rs 0.0", "Rs 0.0", ""),
    @("", "", "", "", "Delivery Subscription", "SM2", "0.0", "Rs 0.0", "Rs 0.0", "Rs 0.0", ""),
    @("", "", "", "", "Stop Delivery Subscription", "SM2", "0.0", "Rs 0.0", "Rs 0.0", "Rs 0.0", ""),
    @("", "", "", "", "Retrive SMS Subscription", "SM2", "0.0", "Rs 0.0", "Rs 0.0", "Rs 0.0", ""),
    @("", "", "", "", "Stop Retrive SMS Subscription", "SM2", "0.0", "Rs 0.0", "Rs 0.0", "Rs 0.0", ""),
    @("", "", "", "", "SMS Inbound Notification", "SM2", "0.0", "Rs 0.0", "Rs 0.0", "Rs 0.0", ""),
    @("", "", "", "", "", "", "Total Amount", "Rs 0.0", "Rs 0.0", "Rs 0.0", "Rs 0.0"),
    @("", "AutInvSP3_AutInvSP3_app3", "smsmessaging", "v1", "Send SMS", "SM1", "0.0", "Rs 0.0", "Rs 0.0", "Rs 0.0", ""),
    @("", "", "", "", "Retrive SMS", "SM2", "0.0", "Rs 0.0", "Rs 0.0", "Rs 0.0", ""),
    @("", "", "", "", "Query SMS Delivery", "SM2", "0.0", "Rs 0.0", "Rs 0.0", "Rs 0.0", ""),
    @("", "", "", "", "Delivery Subscription", "SM2", "0.0", "Rs 0.0", "Rs 0.0", "Rs 0.0", ""),
    @("", "", "", "", "Stop Delivery Subscription", "SM2", "0.0", "Rs 0.0", "Rs 0.0", "Rs 0.0", ""),
    @("", "", "", "", "Retrive SMS Subscription", "SM2", "0.0", "Rs 0.0", "Rs 0.0", "Rs 0.0", ""),
    @("", "", "", "", "Stop Retrive SMS Subscription", "SM2", "0.0", "Rs 0.0", "Rs 0.0", "Rs 0.0", ""),
    @("", "", "", "", "SMS Inbound Notification", "SM2", "0.0", "Rs 0.0", "Rs 0.0", "Rs 0.0", ""),
    @("", "", "", "", "", "", "Total Amount", "Rs 0.0", "Rs 0.0", "Rs 0.0", "Rs 0.0"),
    @("", "", "", "", "", "", "Total Amount", "Rs 0.0", "Rs 0.0", "Rs 0.0", "Rs 0.0"),
    @("", "AutInvSP3_AutInvSP3_app2", "payment", "v1", "Charge", "p1", "0.0", "Rs 0.0", "Rs 0.0", "Rs 0.0", ""),
    @("", "", "", "", "Refund", "RF2", "0.0", "Rs 0.0", "Rs 0.0", "Rs 0.0", ""),
    @("", "", "smsmessaging", "v1", "Send SMS", "SM1", "0.0", "Rs 0.0", "Rs 0.0", "Rs 0.0", ""),
    @("", "", "", "", "Retrive SMS", "SM2", "0.0", "Rs 0.0", "Rs 0.0", "Rs 0.0", ""),
    @("", "", "", "", "Query SMS Delivery", "SM2", "0.0", "Rs 0.0", "Rs 0.0", "Rs 0.0", ""),
    @("", "", "", "", "Delivery Subscription", "SM2", "0.0", "Rs 0.0", "Rs 0.0", "Rs 0.0", ""),
    @("", "", "", "", "Stop Delivery Subscription", "SM2", "0.0", "Rs 0.0", "Rs 0.0", "Rs 0.0", ""),
    @("", "", "", "", "Retrive SMS Subscription", "SM2", "0.0", "Rs 0.0", "Rs 0.0", "Rs 0.0", ""),
    @("", "", "", "", "Stop Retrive SMS Subscription", "SM2", "0.0", "Rs 0.0", "Rs 0.0", "Rs 0.0", ""),
    @("", "", "", "", "SMS Inbound Notification", "SM2", "0.0", "Rs 0.0", "Rs 0.0", "Rs 0.0", ""),
    @("", "", "ussd", "v1", "Send USSD", "u1", "0.0", "Rs 0.0", "Rs 0.0", "Rs 0.0", ""),
    @("", "", "", "", "", "", "Total Amount", "Rs 0.0", "Rs 0.0", "Rs 0.0", "Rs 0.0"),
    @("", "AutInvSP3_AutInvSP3_app1", "payment", "v1", "Charge", "p1", "0.0", "Rs 0.0", "Rs 0.0", "Rs 0.0", ""),
    @("", "", "", "", "Refund", "RF2", "0.0", "Rs 0.0", "Rs 0.0", "Rs 0.0", ""),
    @("", "", "smsmessaging", "v1", "Send SMS", "SM1", "0.0", "Rs 0.0", "Rs 0.0", "Rs 0.0", ""),
    @("", "", "", "", "Retrive SMS", "SM2", "0.0", "Rs 0.0", "Rs 0.0", "Rs 0.0", ""),
    @("", "", "", "", "Query SMS Delivery", "SM2", "0.0", "Rs 0.0", "Rs 0.0", "Rs 0.0", ""),
    @("", "", "", "", "Delivery Subscription", "SM2", "0.0", "Rs 0.0", "Rs 0.0", "Rs 0.0", ""),
    @("", "", "", "", "Stop Delivery Subscription", "SM2", "0.0", "Rs 0.0", "Rs 0.0", "Rs 0.0", ""),
    @("", "", "", "", "Retrive SMS Subscription", "SM2", "0.0", "Rs 0.0", "Rs 0.0", "Rs 0.0", ""),
    @("", "", "", "", "Stop Retrive SMS Subscription", "SM2", "0.0", "Rs 0.0", "Rs 0.0", "Rs 0.0", ""),
    @("", "", "", "", "SMS Inbound Notification", "SM2", "0.0", "Rs 0.0", "Rs 0.0", "Rs 0.0", ""),
    @("", "", "ussd", "v1", "Send USSD", "u1", "0.0", "Rs 0.0", "Rs 0.0", "Rs 0.0", ""),
    @("", "", "", "", "", "", "Total Amount", "Rs 0.0", "Rs 0.0", "Rs 0.0", "Rs 0.0")
)

$r = 1
foreach ($row in $table) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c++
    }
    $r++
}

